# Weekly update: insert a new price record at row 74 for
# "Vega Modelo de Temuco - Chirimoya", pushing the existing rows 74-121
# down to 75-122 (same as Excel's "Insert Row" behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 74, shifting rows 74:121
# down to 75:122.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with this week's record.
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44767
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100107
$ws.Range("H74").Value = "Otros"
$ws.Range("I74").Value = 100107002
$ws.Range("J74").Value = "Chirimoya"
$ws.Range("K74").Value = "Cultivar IV Región"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 75
$ws.Range("N74").Value = 3000
$ws.Range("O74").Value = 3500
$ws.Range("P74").Value = 3300
$ws.Range("Q74").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R74").Value = "Provincia del Elquí"
$ws.Range("S74").Value = 3300
$ws.Range("T74").Value = 1
